$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2: target cluster becomes "FAPs" with updated TPM-derived metrics ---
$ws.Range("D2").Value = "FAPs"
$ws.Range("K2").Value = 2
$ws.Range("L2").Value = 0.6666666666666666
$ws.Range("M2").Value = 0.1817723333333333
$ws.Range("N2").Value = 0.5453170000000001
$ws.Range("O2").Value = 0.008966262009224884
$ws.Range("P2").Value = 0.008966262009224884
$ws.Range("Q2").Value = 0.03935480080066667
$ws.Range("R2").Value = 0.3541932072060001
$ws.Range("S2").Value = 0.008966262009224884
$ws.Range("T2").Value = 0.008966262009224884

# --- Row 3: target cluster becomes "MuSCs" with updated TPM-derived metrics ---
$ws.Range("D3").Value = "MuSCs"
$ws.Range("K3").Value = 3
$ws.Range("L3").Value = 1
$ws.Range("M3").Value = 20.09115
$ws.Range("N3").Value = 60.27345
$ws.Range("O3").Value = 0.9910337379907751
$ws.Range("P3").Value = 0.9910337379907752
$ws.Range("Q3").Value = 4.3498545219
$ws.Range("R3").Value = 39.1486906971
$ws.Range("S3").Value = 0.9910337379907751
$ws.Range("T3").Value = 0.9910337379907752

# --- Row 4: old "MuSCs" row is no longer needed, delete it entirely ---
$ws.Rows.Item(4).Delete()
